$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Old column B ("dbExcel" / Neo4jData.xlsx)
# shifts to C, old column C ("WebExcel" / WebData.xlsx) shifts to D.
$ws.Columns("B").Insert()

# Match new column B's width to column A's width (same visual width as the
# query column).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# New header for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New long Cypher "stat" query value, wrapped the same way as A2.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Lung cancer, NOS']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Move the active selection to the newly added cell.
[void]$ws.Range("B2").Select()

Write-Output "done"
